$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 78, shifting existing rows 78-101 down to 79-102
$ws.Rows.Item(78).Insert()

# Populate the newly inserted row 78 with values (same template as neighboring rows)
$ws.Range("A78").Value = 8
$ws.Range("B78").Value = "Terminal La Palmera de La Serena"
$ws.Range("C78").Value = "Coquimbo"
$ws.Range("D78").Value = 44524
$ws.Range("D78").NumberFormat = $ws.Range("D79").NumberFormat
$ws.Range("E78").Value = 4
$ws.Range("F78").Value = 100112001
$ws.Range("G78").Value = "Berenjena"
$ws.Range("H78").Value = "Sin especificar"
$ws.Range("I78").Value = "Primera"
$ws.Range("J78").Value = 540
$ws.Range("K78").Value = 9000
$ws.Range("L78").Value = 10000
$ws.Range("M78").Value = 9500
$ws.Range("N78").Value = "$/caja 60 unidades"
$ws.Range("O78").Value = "Región de Arica y Parinacota"
$ws.Range("P78").Value = 158
$ws.Range("Q78").Value = 60
$ws.Range("R78").Value = "Hortaliza"
